# Prezenta.xlsx - mark week 8 (column J) attendance as present for several
# students, and move the saved cursor/selection (Newton interpolare, desenare
# grafic).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: J3 was blank, now checked TRUE. Row 3 sits right under the
# header, so its "checked" cells keep the header-adjacent top border (the
# same look as D3/F3/G3 already on this row) - copy that formatting instead
# of the generic "Good" cell style so the border survives.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws.Range("J3").Value = $true

# --- Additional rows where week 8 (column J) becomes checked TRUE.
$checkedRows = @(7, 11, 16, 17, 18, 25)
foreach ($r in $checkedRows) {
    $cell = $ws.Range("J$r")
    $cell.Value = $true
    $cell.Style = "Good"
}

$excel.CutCopyMode = $false

# --- Move the saved view/selection on the sheet.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M12").Select() | Out-Null
